# The "Recorded By" column (G) stored the recorder list as "System, <email>".
# It should instead read "<email>, System" - swap the order of the two
# comma-separated values wherever this exact text occurs in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$colG = $ws.Columns.Item(7)

$first = $colG.Find($oldText)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    $continue = $true
    while ($continue) {
        $current.Value = $newText
        $current = $colG.FindNext($current)
        if ($current -eq $null -or $current.Address() -eq $firstAddress) {
            $continue = $false
        }
    }
}
